# Update RechargeTime (column E) values on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 4
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 2
$ws.Range("E7").Value = 4
$ws.Range("E8").Value = 1
$ws.Range("E9").Value = 3
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 3
$ws.Range("E12").Value = 4
$ws.Range("E13").Value = 3
$ws.Range("E14").Value = 1
$ws.Range("E15").Value = 2
$ws.Range("E16").Value = 3
$ws.Range("E17").Value = 3

# Update the active selection to match the authored view state
$ws.Range("H21").Select()
